# Add a new "Framework" worksheet at the end of the workbook (after "Projects")
# containing two new JavaScript-framework paid courses, matching the commit
# "JavaScript Framework based paid courses added."

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("Projects") so it lands
# at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Framework"

# Column widths matching the other course sheets' layout (B = title, C = url).
$newSheet.Columns.Item(2).ColumnWidth = 72.3
$newSheet.Columns.Item(3).ColumnWidth = 99.3

# Write URL (column C) before title (column B) for each row so new shared
# strings are interned in (url, title) order.
$newSheet.Range("C2").Value = "https://www.udemy.com/course/angular-material-masterclass/"
$newSheet.Range("B2").Value = "Angular Material: Ultimate Masterclass With Angular 9 (2020)"

$newSheet.Range("C4").Value = "https://www.udemy.com/course/complete-react-course-w-hooks-react-router-redux-usecontext/"
$newSheet.Range("B4").Value = "React - The Complete Guide with React Hook Redux 2021 in 4hr"

# Match the selection left on the new active sheet.
$newSheet.Range("C6").Select()
